# Auto-generated Excel COM-interop script to apply the Lich Profits data update
# described by the target diff (scheduled-runner market-price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 900
$ws.Range("J21").Value = 1200
$ws.Range("K21").Value = 900
$ws.Range("L21").Value = 1200
$ws.Range("M21").Value = -432
$ws.Range("N21").Value = -2136

# Row 23
$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 900
$ws.Range("J23").Value = 1200
$ws.Range("K23").Value = 900
$ws.Range("L23").Value = 1200
$ws.Range("M23").Value = -666
$ws.Range("N23").Value = -1668

# Row 38
$ws.Range("H38").Value = 2405.25
$ws.Range("I38").Value = 172.16667
$ws.Range("J38").Value = 9104.5
$ws.Range("K38").Value = 516.50001
$ws.Range("L38").Value = 27313.5
$ws.Range("M38").Value = -144.50001
$ws.Range("N38").Value = -28057.5

# Row 40
$ws.Range("H40").Value = 27277450
$ws.Range("I40").Value = 7497
$ws.Range("J40").Value = 42860280
$ws.Range("K40").Value = 7497
$ws.Range("L40").Value = 42860280
$ws.Range("M40").Value = -7322
$ws.Range("N40").Value = -42860630

# Row 51
$ws.Range("H51").Value = 13360.516
$ws.Range("I51").Value = 17664.834
$ws.Range("J51").Value = 12404
$ws.Range("K51").Value = 17664.834
$ws.Range("L51").Value = 12404
$ws.Range("M51").Value = -17180.834
$ws.Range("N51").Value = -13372

# Row 58
$ws.Range("H58").Value = 409.8125
$ws.Range("I58").Value = 409.8125
$ws.Range("K58").Value = 1229.4375
$ws.Range("M58").Value = -1079.4375

# Row 64
$ws.Range("H64").Value = 3516.2
$ws.Range("I64").Value = 3562
$ws.Range("J64").Value = 3333
$ws.Range("K64").Value = 3562
$ws.Range("L64").Value = 3333
$ws.Range("M64").Value = -3314
$ws.Range("N64").Value = -3829

# Row 67
$ws.Range("H67").Value = 3516.2
$ws.Range("I67").Value = 3562
$ws.Range("J67").Value = 3333
$ws.Range("K67").Value = 3562
$ws.Range("L67").Value = 3333
$ws.Range("M67").Value = -2704
$ws.Range("N67").Value = -5049

# Row 74
$ws.Range("H74").Value = 4440.2856
$ws.Range("I74").Value = 4620.4
$ws.Range("K74").Value = 4620.4
$ws.Range("M74").Value = -3684.4

# Row 76
$ws.Range("H76").Value = 4225.857
$ws.Range("J76").Value = 4789.3335
$ws.Range("L76").Value = 4789.3335
$ws.Range("N76").Value = -5419.3335

# Row 77
$ws.Range("H77").Value = 4440.2856
$ws.Range("I77").Value = 4620.4
$ws.Range("K77").Value = 23102
$ws.Range("M77").Value = -18422

# Row 79
$ws.Range("H79").Value = 4225.857
$ws.Range("J79").Value = 4789.3335
$ws.Range("L79").Value = 4789.3335
$ws.Range("N79").Value = -6973.3335

# Row 100
$ws.Range("H100").Value = 1656.5
$ws.Range("I100").Value = 1580.7142
$ws.Range("K100").Value = 1580.7142
$ws.Range("M100").Value = -1039.7142

# Row 132
$ws.Range("H132").Value = 3165.5305
$ws.Range("I132").Value = 2636.425
$ws.Range("K132").Value = 7909.275000000001
$ws.Range("M132").Value = -5379.275000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 34
$ws.Range("H34").Value = 49994
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 49994
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 49994
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = -50536

# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = ""

# Row 97
$ws.Range("H97").Value = 827.9474
$ws.Range("I97").Value = 674.7646999999999
$ws.Range("K97").Value = 674.7646999999999
$ws.Range("M97").Value = -178.7646999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 747.9231
$ws.Range("I94").Value = 585.1429000000001
$ws.Range("J94").Value = 1431.6
$ws.Range("K94").Value = 585.1429000000001
$ws.Range("L94").Value = 1431.6
$ws.Range("M94").Value = -134.1429000000001
$ws.Range("N94").Value = -2333.6

# Row 99
$ws.Range("H99").Value = 5391.2104
$ws.Range("I99").Value = 8766.666999999999
$ws.Range("J99").Value = 4758.3125
$ws.Range("K99").Value = 8766.666999999999
$ws.Range("L99").Value = 4758.3125
$ws.Range("M99").Value = -7268.666999999999
$ws.Range("N99").Value = -7754.3125

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 314363.2
$ws.Range("I31").Value = 372023.12
$ws.Range("K31").Value = 372023.12
$ws.Range("M31").Value = -371728.12

# Row 34
$ws.Range("H34").Value = 314363.2
$ws.Range("I34").Value = 372023.12
$ws.Range("K34").Value = 372023.12
$ws.Range("M34").Value = -371821.12

$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 87.666664
$ws.Range("I6").Value = 105.666664
$ws.Range("J6").Value = 69.666664
$ws.Range("K6").Value = 316.999992
$ws.Range("L6").Value = 208.999992
$ws.Range("M6").Value = -203.999992
$ws.Range("N6").Value = -434.999992

# Row 12
$ws.Range("H12").Value = 108.4
$ws.Range("J12").Value = 148.6
$ws.Range("L12").Value = 445.8
$ws.Range("N12").Value = -791.8

# Row 70
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = ""

# Row 73
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = ""

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 10000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = ""

# Row 92
$ws.Range("H92").Value = 12042.714
$ws.Range("J92").Value = 12042.714
$ws.Range("L92").Value = 12042.714
$ws.Range("N92").Value = -15786.714

$ws = $wb.Worksheets.Item("LTW")
# Row 38
$ws.Range("H38").Value = 79332.664
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 79332.664
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 79332.664
$ws.Range("M38").Value = ""
$ws.Range("N38").Value = -80152.664

# Row 41
$ws.Range("H41").Value = 49999
$ws.Range("I41").Value = 49999
$ws.Range("K41").Value = 49999
$ws.Range("M41").Value = -49561

# Row 50
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").Value = ""

$ws = $wb.Worksheets.Item("WVR")
# Row 34
$ws.Range("H34").Value = 10025
$ws.Range("I34").Value = 10025
$ws.Range("K34").Value = 10025
$ws.Range("M34").Value = -9822

# Row 37
$ws.Range("H37").Value = 65162.6
$ws.Range("J37").Value = 68953.75
$ws.Range("L37").Value = 68953.75
$ws.Range("N37").Value = -69359.75

# Row 40
$ws.Range("H40").Value = 37999
$ws.Range("I40").Value = 37999
$ws.Range("K40").Value = 37999
$ws.Range("M40").Value = -37850

# Row 42
$ws.Range("H42").Value = 49997
$ws.Range("I42").Value = 49997
$ws.Range("K42").Value = 49997
$ws.Range("M42").Value = -49619

# Row 43
$ws.Range("H43").Value = 39998
$ws.Range("I43").Value = 39998
$ws.Range("K43").Value = 39998
$ws.Range("M43").Value = -39849

# Row 126
$ws.Range("H126").Value = 2258.95
$ws.Range("I126").Value = 2149.625
$ws.Range("K126").Value = 6448.875
$ws.Range("M126").Value = -3978.875

# Row 136
$ws.Range("H136").Value = 701341.7
$ws.Range("I136").Value = 834677.2
$ws.Range("K136").Value = 2504031.6
$ws.Range("M136").Value = -2501481.6
